# Insert a new weekly record at row 208, pushing the existing rows
# 208-218 down to 209-219 (the dataset grows from A1:R218 to A1:R219).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 208, shifting rows
# 208:218 down to 209:219.
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new weekly observation.
$ws.Cells.Item(208, 1).Value = 4
$ws.Cells.Item(208, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(208, 3).Value = "Los Lagos"
$ws.Cells.Item(208, 4).Value = 44610
$ws.Cells.Item(208, 5).Value = 10
$ws.Cells.Item(208, 6).Value = 100112021
$ws.Cells.Item(208, 7).Value = "Ají"
$ws.Cells.Item(208, 8).Value = "Inferno"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 130
$ws.Cells.Item(208, 11).Value = 18000
$ws.Cells.Item(208, 12).Value = 18000
$ws.Cells.Item(208, 13).Value = 18000
$ws.Cells.Item(208, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(208, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(208, 16).Value = 1500
$ws.Cells.Item(208, 17).Value = 12
$ws.Cells.Item(208, 18).Value = "Hortaliza"
